$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-order the "Aged Care Active Outbreaks" metric rows (C65:C71) ---
# Current order (rows 65-71): Facilities(Weekly), Facilities(Weekly) per 1M,
#   Facilities(Weekly) Change, Active Outbreaks, (7-day avg), (7-day avg) per 1M,
#   Weekly Change
# Target order: Active Outbreaks, (7-day avg), (7-day avg) per 1M, Weekly Change,
#   Facilities(Weekly), Facilities(Weekly) per 1M, Facilities(Weekly) Change
# i.e. rotate the 7 cells left by 3 (move the first 3 rows to the bottom),
# carrying their values *and* formatting (style s="2") along with them.

# Stash the first three rows' C-column cells (with formatting) in a scratch
# area well outside the used range.
$ws.Range("C65:C67").Copy($ws.Range("Z65:Z67"))

# Slide the remaining four rows up into C65:C68.
$ws.Range("C68:C71").Copy($ws.Range("C65:C68"))

# Drop the stashed three rows back in at the bottom (C69:C71).
$ws.Range("Z65:Z67").Copy($ws.Range("C69:C71"))

# Clean up the scratch cells.
$ws.Range("Z65:Z67").Clear()

# --- 2. Rename the "Active Outbreaks Facilities" metrics (drop "Active") ---
for ($r = 69; $r -le 71; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2.Replace("Active Outbreaks Facilities", "Outbreaks Facilities")
}

# --- 3. Update the saved selection to match the author's final cursor spot ---
$ws.Range("D70").Select()
